# Update the cryptocurrency price / link / volume table on Sheet1.
# Generated from the target diff: refreshed Price (D) and Volume(1h) (E)
# figures for most rows, plus three pairs of rows that swapped rank
# (12<->13, 22<->23, 50<->51), each carrying its Coin name (B), Link (C),
# Price (D) and Volume (E) along with it.
#
# Cells in column D whose new text happens to look like a plain number
# (e.g. "238.12") are written with a leading apostrophe so Excel keeps
# storing them as text (matching the original inlineStr cells) instead of
# silently re-typing them as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.254.80"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.856.55"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'0.7031"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").Value = "'238.12"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.08010"
$ws.Range("E8").Value = "  +5.31%  "
$ws.Range("D9").Value = "'0.3030"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "'23.60"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "'0.08196"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.849.91"
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.201"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'0.7072"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "'89.77"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "29.194.35"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'5.836"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").Value = "'0.000007853"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "'13.23"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'236.95"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.075.17"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "'7.517"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").Value = "'163.18"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "'8.896"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "'18.12"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "'1.911"
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("D30").Value = "'1.402"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "'1.475"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "'4.348"
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("D33").Value = "'4.030"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").Value = "'0.05191"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "'1.167"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").Value = "'0.7157"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").Value = "'0.9953"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'2.722"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "'0.9375"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").Value = "1.148.92"
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").Value = "'5.993"
$ws.Range("D44").Value = "'0.4263"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "'70.37"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'103.00"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "'0.5287"
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("D49").Value = "'1.744"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "1.981.22"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.175"
$ws.Range("E51").Value = "  +0.42%  "
